$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, pushing old rows 15-20 down to 16-21.
# (Excel automatically adjusts the total row's SUM/QUOTIENT/MOD formula ranges.)
$ws.Rows(15).Insert()

# Update row 13: was "1/3-2018" / "Enemy implementation" / 0h30m
# becomes a new log entry: "28/2-2018" / "Game implementation" / 1h30m
$ws.Range("A13").Value = "28/2-2018"
$ws.Range("B13").Value = "Game implementation"
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 30

# Row 14 keeps "1/3-2018" but now holds the "Enemy implementation" entry
# that used to live in row 13 (0h30m, unchanged).
$ws.Range("A14").Value = "1/3-2018"
$ws.Range("B14").Value = "Enemy implementation"
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 30

# New row 15 (freshly inserted): "1/3-2018" / "Player basic completed" / 0h30m
$ws.Range("A15").Value = "1/3-2018"
$ws.Range("B15").Value = "Player basic completed"
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 30

# Row 16 already holds the old row 15 contents ("1/3-2018" / "Game completed,
# transformable and overloading" / 3h0m) thanks to the row insert shifting it down.

# Update the sheet view's selected cell to match the author's final position.
$ws.Range("C19").Select()
